$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12 used to hold the "Test unit 07" / "Rosepetale" result (with no Date value).
# It is repurposed to hold the "Test unit 06" / "Rosepetale" result with an updated Date
# and new measurements.
$ws.Range("A12").Value = "Test unit 06"
$ws.Range("B12").Value = "August 12th-14th"
$ws.Range("C12").Value = 16
$ws.Range("E12").Value = 20

# Row 13 is a brand-new data row that now carries the "Test unit 07" / "Rosepetale"
# result that used to live in row 12. Copy row 12's formatting first so the new row
# matches the look of the other data rows, then fill in its own values/formulas.
$ws.Range("A12:J12").Copy() | Out-Null
$ws.Range("A13:J13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A13").Value = "Test unit 07"
$ws.Range("B13").Value = "August 12th-14th"
$ws.Range("C13").Value = 16
$ws.Range("D13").Value = "Rosepetale"
$ws.Range("E13").Value = 20.5
$ws.Range("F13").Value = 50
$ws.Range("G13").Formula = "=(F13-E13)"
$ws.Range("H13").Formula = "=(G13/C13)"
$ws.Range("I13").Formula = "=(H13/3600)"
$ws.Range("J13").Value = "400 Sec"

# Fix the typo in the "Date" column for the Sensual musk rows: "Augus 7th - 12th" -> "August 7th - 12th"
$ws.Range("B4").Value = "August 7th - 12th"
$ws.Range("B5").Value = "August 7th - 12th"
$ws.Range("B6").Value = "August 7th - 12th"

# Update the active cell selection to match the saved view state.
$ws.Range("F8").Select() | Out-Null

$wb.Save()
